$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tf = $s.Shapes.Item(1).TextFrame
$tr = $tf.TextRange

# The "small caps" paragraph is the 2nd paragraph in the text box.
$smallCapsPara = $tr.Paragraphs(2, 1)

# Insert a brand-new paragraph right after it (before the "subscripts" paragraph),
# seeded with the plain text first - formatting is applied afterwards.
$smallCapsPara.InsertAfter("`rHere is some underlined text") | Out-Null

# Grab the paragraph we just created (now paragraph #3).
$newPara = $tr.Paragraphs(3, 1)

# "Here is " -> left as-is (no special formatting)
# "some "    -> underline
$newPara.Characters(9, 5).Font.Underline = $true
# "underlined" -> italic + underline
$r = $newPara.Characters(14, 10)
$r.Font.Italic = $true
$r.Font.Underline = $true
# " " -> underline
$newPara.Characters(24, 1).Font.Underline = $true
# "text" -> bold + underline
$r = $newPara.Characters(25, 4)
$r.Font.Bold = $true
$r.Font.Underline = $true
